$wb = $excel.ActiveWorkbook

# --- Sheet "Summary" ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B4").Value = 78
$ws1.Range("B5").Value = 22
$ws1.Range("B6").Value = "Amber: Fee Drag > 10%"
$ws1.Range("B7").Value = 22

# --- Sheet "Symbols" ---
$ws2 = $wb.Worksheets.Item("Symbols")
$ws2.Range("B2").Value = 78

# --- Sheet "Strategies" ---
$ws3 = $wb.Worksheets.Item("Strategies")
$ws3.Range("D2").Value = 78
$ws3.Range("E2").Value = 100
$ws3.Range("F2").Value = 22
$ws3.Range("G2").Value = 4
$ws3.Range("H2").Value = 19.5
$ws3.Range("J2").Value = "[{'strategy_name': 'Short Put', 'pnl': 100.0, 'entry_ts': '2025-01-01T10:00:00', 'exit_ts': '2025-01-05T10:00:00'}]"
